$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update movie/video titles in column A (shortened / simplified titles) ---
# Order matches the order these new strings were introduced by the original edit.
$ws.Range("A2").Value  = 'iron man vs loki'
$ws.Range("A3").Value  = 'avengers age of ultron'
$ws.Range("A5").Value  = "the devil wears prada andy's interview"
$ws.Range("A4").Value  = 'the devil wears prada andy gets a makeover'
$ws.Range("A7").Value  = 'inception ending'
$ws.Range("A6").Value  = 'inception the escape from limbo'
$ws.Range("A10").Value = 'no time to die'
$ws.Range("A9").Value  = 'hot n cold'
$ws.Range("A8").Value  = 'i love me'

# --- Update the "Louis Poulsen" product text (row 4) to drop the accented French wording ---
$ws.Range("F4").Value = 'Louis Poulsen.png, blazer with patch detail.png'
$ws.Range("C4").Value = 'Louis Poulsen, blazer with patch detail'

# --- Resize columns to fit the new data layout ---
$ws.Columns.Item(1).ColumnWidth = 63 + 1/4
$ws.Columns.Item(2).ColumnWidth = 37 + 3/4
$ws.Columns.Item(3).ColumnWidth = 38 + 1/12
$ws.Columns.Item(5).ColumnWidth = 15 + 7/12

# --- Move the active selection ---
$ws.Range("C8").Select()
